# Script 1 - atualização em 2025-09-20 17:07:42Z
# Applies updated values to the "Valor" (D) column and one "Posição" (E) cell
# on the single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value  = 1.305594974582545
$ws.Range("D7").Value  = 1.016936609646641
$ws.Range("D8").Value  = 0.9260032757814711
$ws.Range("D9").Value  = 0.7711935476200089
$ws.Range("D10").Value = 0.5859108634114807
$ws.Range("D11").Value = 0.4750618247492207
$ws.Range("D12").Value = 0.2494238130431646
$ws.Range("D13").Value = 1.469372305885356
$ws.Range("D19").Value = 0.9747549918484564
$ws.Range("D20").Value = 0.7955823415593115
$ws.Range("D21").Value = 0.5829578958620336
$ws.Range("D23").Value = 0.2724902505250656
$ws.Range("E34").Value = 25
